$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UI Test")

# --- Row 5 (UI-02 test case): the Expected Results cell is updated ---
$ws.Range("L5").Value = "Validation message appears for each incorrect field"

# --- Row 6: new UI-03 test case ("Verify system responsiveness on mobile view") ---
$ws.Range("C6").Value = "UI-03"
$ws.Range("D6").Value = "Verify system display is responsive on mobile devices"
$ws.Range("E6").Value = "UI-TC-03"
$ws.Range("F6").Value = "UI Module"
$ws.Range("G6").Value = "UI-TS-03"
$ws.Range("H6").Value = "Verify system responsiveness on mobile view"
$ws.Range("I6").Value = "Application is accessible"
$ws.Range("J6").Value = "1. Open application in browser " + [char]10 + "2. Switch to mobile view (responsive mode) " + [char]10 + "3. Navigate through pages"
$ws.Range("K6").Value = "Mobile screen resolution"
$ws.Range("L6").Value = "Layout adjusts properly, no overlapping text/buttons, all features accessible"

# --- View state: move the selection to match the saved workbook state ---
$ws.Activate()
$ws.Range("L7").Select()
